$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 293, shifting existing rows 293:334 down to 294:335.
$ws.Rows.Item(293).Insert()

# Populate the newly inserted row 293 with the new weekly price observation.
$ws.Cells.Item(293, 1).Value = 7
$ws.Cells.Item(293, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(293, 3).Value = "Ñuble"
$ws.Cells.Item(293, 4).Value = 45131
$ws.Cells.Item(293, 5).Value = 16
$ws.Cells.Item(293, 6).Value = 100112032
$ws.Cells.Item(293, 7).Value = "Zapallo italiano"
$ws.Cells.Item(293, 8).Value = "Sin especificar"
$ws.Cells.Item(293, 9).Value = "Primera"
$ws.Cells.Item(293, 10).Value = 100
$ws.Cells.Item(293, 11).Value = 14000
$ws.Cells.Item(293, 12).Value = 14000
$ws.Cells.Item(293, 13).Value = 14000
$ws.Cells.Item(293, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(293, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(293, 16).Value = 280
$ws.Cells.Item(293, 17).Value = 50
$ws.Cells.Item(293, 18).Value = "Hortaliza"
